$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 45975
$ws.Range("D8").Value = 172.24
$ws.Range("E8").Value = 161.46
$ws.Range("F8").Value = 171.46
$ws.Range("G8").Value = 161.63

# Row 9
$ws.Range("A9").Value = 45975
$ws.Range("D9").Value = 172.24
$ws.Range("E9").Value = 161.46
$ws.Range("F9").Value = 171.46
$ws.Range("G9").Value = 161.63

# Row 10
$ws.Range("A10").Value = 45975
$ws.Range("D10").Value = 174.8
$ws.Range("E10").Value = 164.66
$ws.Range("F10").Value = 174.66
$ws.Range("G10").Value = 165.16

# Row 11
$ws.Range("A11").Value = 45974
$ws.Range("D11").Value = 171.51
$ws.Range("E11").Value = 160.94
$ws.Range("F11").Value = 170.94
$ws.Range("G11").Value = 161.1

# Row 12
$ws.Range("A12").Value = 45974
$ws.Range("D12").Value = 171.51
$ws.Range("E12").Value = 160.94
$ws.Range("F12").Value = 170.94
$ws.Range("G12").Value = 161.1

# Row 13
$ws.Range("A13").Value = 45974
$ws.Range("D13").Value = 173.97
$ws.Range("E13").Value = 163.94
$ws.Range("F13").Value = 173.94
$ws.Range("G13").Value = 164.43

# Row 17
$ws.Range("A17").Value = 45975
$ws.Range("D17").Value = 178.16
$ws.Range("E17").Value = 167.49
$ws.Range("F17").Value = 177.49

# Row 18
$ws.Range("A18").Value = 45974
$ws.Range("D18").Value = 177.35
$ws.Range("E18").Value = 166.8
$ws.Range("F18").Value = 176.8

# Row 22
$ws.Range("A22").Value = 45975
$ws.Range("D22").Value = 173.59
$ws.Range("E22").Value = 163.74
$ws.Range("F22").Value = 173.34
$ws.Range("G22").Value = 165.03

# Row 23
$ws.Range("A23").Value = 45975
$ws.Range("D23").Value = 179.58
$ws.Range("E23").Value = 168.27
$ws.Range("F23").Value = 178.27

# Row 24
$ws.Range("A24").Value = 45975
$ws.Range("D24").Value = 179.38
$ws.Range("E24").Value = 168.5
$ws.Range("F24").Value = 178.5

# Row 25
$ws.Range("A25").Value = 45975
$ws.Range("D25").Value = 180.21
$ws.Range("E25").Value = 167.91
$ws.Range("F25").Value = 177.91
$ws.Range("G25").Value = 167.95

# Row 26
$ws.Range("A26").Value = 45975
$ws.Range("D26").Value = 178.91
$ws.Range("E26").Value = 169.48
$ws.Range("F26").Value = 179.48

# Row 27
$ws.Range("A27").Value = 45974
$ws.Range("D27").Value = 172.75
$ws.Range("E27").Value = 163.03
$ws.Range("F27").Value = 172.63
$ws.Range("G27").Value = 164.31

# Row 28
$ws.Range("A28").Value = 45974
$ws.Range("D28").Value = 178.75
$ws.Range("E28").Value = 167.66
$ws.Range("F28").Value = 177.66

# Row 29
$ws.Range("A29").Value = 45974
$ws.Range("D29").Value = 178.55
$ws.Range("E29").Value = 167.88
$ws.Range("F29").Value = 177.88

# Row 30
$ws.Range("A30").Value = 45974
$ws.Range("D30").Value = 179.38
$ws.Range("E30").Value = 167.28
$ws.Range("F30").Value = 177.28
$ws.Range("G30").Value = 167.32

# Row 31
$ws.Range("A31").Value = 45974
$ws.Range("D31").Value = 178.09
$ws.Range("E31").Value = 168.85
$ws.Range("F31").Value = 178.85

# Row 35
$ws.Range("A35").Value = 45975
$ws.Range("D35").Value = 173.16
$ws.Range("E35").Value = 161.76
$ws.Range("F35").Value = 170.76

# Row 36
$ws.Range("A36").Value = 45974
$ws.Range("D36").Value = 172.44
$ws.Range("E36").Value = 161.15
$ws.Range("F36").Value = 170.15

# Row 40
$ws.Range("A40").Value = 45975
$ws.Range("D40").Value = 178.77
$ws.Range("E40").Value = 167.35
$ws.Range("F40").Value = 177.35

# Row 41
$ws.Range("A41").Value = 45975
$ws.Range("D41").Value = 178.48
$ws.Range("E41").Value = 167.77
$ws.Range("F41").Value = 177.77

# Row 42
$ws.Range("A42").Value = 45974
$ws.Range("D42").Value = 177.92
$ws.Range("E42").Value = 166.63
$ws.Range("F42").Value = 176.63

# Row 43
$ws.Range("A43").Value = 45974
$ws.Range("D43").Value = 177.63
$ws.Range("E43").Value = 167.05
$ws.Range("F43").Value = 177.05

# Row 47
$ws.Range("A47").Value = 45975
$ws.Range("D47").Value = 172.16
$ws.Range("E47").Value = 162.59
$ws.Range("F47").Value = 172.59

# Row 48
$ws.Range("A48").Value = 45975
$ws.Range("D48").Value = 172.14
$ws.Range("E48").Value = 162.75
$ws.Range("F48").Value = 172.75

# Row 49
$ws.Range("A49").Value = 45974
$ws.Range("D49").Value = 171.39
$ws.Range("E49").Value = 162.29
$ws.Range("F49").Value = 172.29

# Row 50
$ws.Range("A50").Value = 45974
$ws.Range("D50").Value = 171.38
$ws.Range("E50").Value = 162.46
$ws.Range("F50").Value = 172.46

# Row 54
$ws.Range("A54").Value = 45975
$ws.Range("D54").Value = 188.88
$ws.Range("E54").Value = 177.86
$ws.Range("F54").Value = 187.86

# Row 55
$ws.Range("A55").Value = 45975
$ws.Range("D55").Value = 176.56
$ws.Range("E55").Value = 174.95
$ws.Range("F55").Value = 184.95

# Row 56
$ws.Range("A56").Value = 45975
$ws.Range("D56").Value = 179.06

# Row 57
$ws.Range("A57").Value = 45975
$ws.Range("D57").Value = 178.59
$ws.Range("E57").Value = 169.22

# Row 58
$ws.Range("A58").Value = 45975
$ws.Range("D58").Value = 174.5
$ws.Range("E58").Value = 165.27
$ws.Range("F58").Value = 175.27

# Row 59
$ws.Range("A59").Value = 45975
$ws.Range("D59").Value = 181.13
$ws.Range("E59").Value = 175.9

# Row 60
$ws.Range("A60").Value = 45974
$ws.Range("D60").Value = 188.06
$ws.Range("E60").Value = 177.14
$ws.Range("F60").Value = 187.14

# Row 61
$ws.Range("A61").Value = 45974
$ws.Range("D61").Value = 175.73
$ws.Range("E61").Value = 174.36
$ws.Range("F61").Value = 184.36

# Row 62
$ws.Range("A62").Value = 45974
$ws.Range("D62").Value = 178.22

# Row 63
$ws.Range("A63").Value = 45974
$ws.Range("D63").Value = 177.79
$ws.Range("E63").Value = 168.63

# Row 64
$ws.Range("A64").Value = 45974
$ws.Range("D64").Value = 173.7
$ws.Range("E64").Value = 164.68
$ws.Range("F64").Value = 174.68

# Row 65
$ws.Range("A65").Value = 45974
$ws.Range("D65").Value = 180.33
$ws.Range("E65").Value = 175.24
